$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: S.No=3, Date=45612 (2024-11-16), DSA/OS/DSA_Interview/DEV_Interview = 0hr00min, Project/Total = 4hr00min
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 45612
$ws.Range("F5").Value = "0hr00min"
$ws.Range("G5").Value = "0hr00min"
$ws.Range("H5").Value = "4hr00min"
$ws.Range("I5").Value = "0hr00min"
$ws.Range("J5").Value = "0hr00min"
$ws.Range("K5").Value = "4hr00min"

# Row 6: S.No=4, Date=45247 (2023-11-17), DSA/OS/DSA_Interview/DEV_Interview = 0hr00min, Project/Total = 6hr00min
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 45247
$ws.Range("F6").Value = "0hr00min"
$ws.Range("G6").Value = "0hr00min"
$ws.Range("H6").Value = "6hr00min"
$ws.Range("I6").Value = "0hr00min"
$ws.Range("J6").Value = "0hr00min"
$ws.Range("K6").Value = "6hr00min"

# Copy the date formatting (style) from an existing date cell onto the new date cells
$ws.Range("E3").Copy()
$ws.Range("E5:E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selection to match the final active cell
$ws.Range("K6").Select()
